$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-8 from 45243 to 45244 (date serial +1 day)
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45244
}
